$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first_eval")

$ws.Range("B2").Value = 0.3934793239648547
$ws.Range("C2").Value = 0.4037288085692023
$ws.Range("D2").Value = 0.2175455044138317
$ws.Range("E2").Value = 0.4664177359554756
$ws.Range("F2").Value = 0.2626622900771862
$ws.Range("G2").Value = 11
$ws.Range("B3").Value = 0.3075541611523277
$ws.Range("C3").Value = 0.3490093155375321
$ws.Range("D3").Value = 0.1614430495769377
$ws.Range("E3").Value = 0.4017997630374335
$ws.Range("F3").Value = 0.2691182110078417
$ws.Range("G3").Value = 13
$ws.Range("B4").Value = 0.2355384709971136
$ws.Range("C4").Value = 0.2831383568828125
$ws.Range("D4").Value = 0.1164317108812595
$ws.Range("E4").Value = 0.3412209121394226
$ws.Range("F4").Value = 0.2578653762121229
$ws.Range("G4").Value = 12
$ws.Range("B5").Value = 0.3198529258979658
$ws.Range("C5").Value = 0.3395672504785434
$ws.Range("D5").Value = 0.1389581556796866
$ws.Range("E5").Value = 0.3727709158178606
$ws.Range("F5").Value = 0.2030610601727263
$ws.Range("G5").Value = 9
$ws.Range("B6").Value = 0.3074279757554496
$ws.Range("C6").Value = 0.3343017199468005
$ws.Range("D6").Value = 0.1343402168544869
$ws.Range("E6").Value = 0.3665245105780607
$ws.Range("F6").Value = 0.2116761409549216
$ws.Range("G6").Value = 9
$ws.Range("B7").Value = 0.2926088612260327
$ws.Range("C7").Value = 0.3274254659629177
$ws.Range("D7").Value = 0.1289431368268439
$ws.Range("E7").Value = 0.3590865311131063
$ws.Range("F7").Value = 0.2207681817058433
$ws.Range("G7").Value = 9
$ws.Range("B8").Value = 0.317927339586861
$ws.Range("C8").Value = 0.3473839361765009
$ws.Range("D8").Value = 0.1472957456731997
$ws.Range("E8").Value = 0.3837912787873113
$ws.Range("F8").Value = 0.2298271460192152
$ws.Range("G8").Value = 8
$ws.Range("B9").Value = 0.3052654375636107
$ws.Range("C9").Value = 0.3373465939268102
$ws.Range("D9").Value = 0.1413070844341813
$ws.Range("E9").Value = 0.3759083457894774
$ws.Range("F9").Value = 0.23693904963195
$ws.Range("G9").Value = 7
$ws.Range("B10").Value = 0.3173063146244254
$ws.Range("C10").Value = 0.3418763137882747
$ws.Range("D10").Value = 0.1463084296434522
$ws.Range("E10").Value = 0.3825028491965154
$ws.Range("F10").Value = 0.2307148046110977
$ws.Range("G10").Value = 7
$ws.Range("B11").Value = 0.3606037648954714
$ws.Range("C11").Value = 0.379698193710172
$ws.Range("D11").Value = 0.1799938953220499
$ws.Range("E11").Value = 0.4242568742189688
$ws.Range("F11").Value = 0.2448480836729456
$ws.Range("G11").Value = 6
